# Fully working Calibration - Four Buttons + Annotations (slate) Version
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worlds-like test")

# --- Update calibration data values (rows 2-5, "A1, Pointing at ID2") ---
$ws.Range("B2").Value = -49.3
$ws.Range("C2").Value = -47.7

$ws.Range("B3").Value = 23.5
$ws.Range("C3").Value = -48.2

$ws.Range("B4").Value = 47.9
$ws.Range("C4").Value = -23.6

$ws.Range("C5").Value = -23.3

# --- Update calibration data values (rows 30-33, second data block) ---
$ws.Range("B30").Value = -23
$ws.Range("C30").Value = -47.4

$ws.Range("B31").Value = -47.1
$ws.Range("C31").Value = -47.6

$ws.Range("B32").Value = -24.5
$ws.Range("C32").Value = -23

$ws.Range("B33").Value = -46.7
$ws.Range("C33").Value = -23.2

# --- Update view/selection state ---
$ws.Activate()
$ws.Range("K2:K5").Select()
